# Small board "Ping" pin assignment: prepend "PA6/" as its own run in
# front of the existing "CCP1" run, so the paragraph reads "PA6/CCP1"
# (matching the diff's two-run structure rather than merging into one).
$d = $word.ActiveDocument

# Locate the paragraph that currently holds just "CCP1" (directly under
# the "Ping" paragraph, small-board section).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.Trim() -eq "CCP1") {
        $target = $para
        break
    }
}

$r = $target.Range
$r.Collapse(1)                      # wdCollapseStart

# Insert a placeholder run, give it different formatting so the engine
# keeps it as a distinct <w:r> instead of coalescing it into the
# following "CCP1" run, then fix up the text and formatting afterwards
# so the final run matches the original run's properties exactly.
$r.InsertBefore("PA6X")
$newRun = $d.Range($target.Range.Start, $target.Range.Start + 4)
$newRun.Font.Size = 20

$find = $target.Range.Find
$find.ClearFormatting()
$find.Execute("PA6X", $true, $false, $false, $false, $false, $true, 1, $false, "PA6/", 2)

$newRun2 = $d.Range($target.Range.Start, $target.Range.Start + 4)
$newRun2.Font.Size = 12
